# LOQ4244.xlsx — course-outline sheet rewritten.
#
# The course description ("Lean Six Sigma") was replaced with a generic
# placeholder course record: most of the long Portuguese/English narrative
# cells (objectives, long syllabus, bibliography) were removed, a couple of
# short header/value cells were inserted ("Critério:" / the professor line
# moved up under "Objetivos:"), and the sheet shrank from 24 to 23 used rows.
# We rewrite the whole A1:C23 block explicitly rather than trying to
# reverse-engineer individual row inserts/deletes, then fix up the handful
# of row heights that changed and drop the now-unused row 24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @(
    @{ Row = 1; A = $null; B = 'Ementa atual:'; C = 'Ementa modificada (dados modificados em vermelho):' },
    @{ Row = 2; A = $null; B = 'LOQ4244'; C = 'LOQ4244' },
    @{ Row = 3; A = 'Nome:'; B = ' Lean Six Sigma'; C = ' Lean Six Sigma' },
    @{ Row = 4; A = 'Name:'; B = 'Lean Six Sigma'; C = 'Lean Six Sigma' },
    @{ Row = 5; A = 'Créditos-aula:'; B = '2'; C = '2' },
    @{ Row = 6; A = 'Créditos-trabalho'; B = '1'; C = '1' },
    @{ Row = 7; A = 'Carga horária:'; B = '60 h'; C = '60 h' },
    @{ Row = 8; A = 'Ativação:'; B = '01/01/2018'; C = '01/01/2018' },
    @{ Row = 9; A = 'Semestre ideal:'; B = 'EP-9'; C = 'EP-9' },
    @{ Row = 10; A = 'Objetivos:'; B = '5840535 - Messias Borges Silva'; C = '5840535 - Messias Borges Silva' },
    @{ Row = 11; A = 'Objectives:'; B = 'Provide basic knowledge and applications of the usual techniques of lean manufacturing and six sigma, presenting concepts, principles and tools used in this approach. Emphasis will be placed on lean thinking, value stream mapping, lean production control systems design, manufacturing cell design, multifunction worker development, and requirements definition for the implementation of lean production lines'; C = 'Provide basic knowledge and applications of the usual techniques of lean manufacturing and six sigma, presenting concepts, principles and tools used in this approach. Emphasis will be placed on lean thinking, value stream mapping, lean production control systems design, manufacturing cell design, multifunction worker development, and requirements definition for the implementation of lean production lines' },
    @{ Row = 12; A = 'Docentes responsáveis:'; B = $null; C = $null },
    @{ Row = 13; A = 'Programa resumido:'; B = 'Semestral'; C = 'Semestral' },
    @{ Row = 14; A = 'Short syllabus:'; B = 'Lean thinking; Waste Categories; Understanding the Value Stream; Value Stream Mapping; Takt Time; Production Pull System ; Creation of continuous flow ; Manufacturing Cells; Elements of Production Control; Tools for Lean production .Six Sigma Strategy.'; C = 'Lean thinking; Waste Categories; Understanding the Value Stream; Value Stream Mapping; Takt Time; Production Pull System ; Creation of continuous flow ; Manufacturing Cells; Elements of Production Control; Tools for Lean production .Six Sigma Strategy.' },
    @{ Row = 15; A = 'Programa:'; B = '01/01/2018'; C = '01/01/2018' },
    @{ Row = 16; A = 'Syllabus:'; B = 'The Lean Production Approach; History; Lean Thonking; Waste Categories; DMAIC methodology applied to the Lean Six Sigma; Organizational culture for the Lean; Understanding the Value Stream; Value Stream Mapping; Flow of Material and Information; Characteristics of the Lean Value Stream; Concept of Takt Time; Capability Analysis; Production System Pulled; Creation of continuous flow of production; Manufacturing Cells; Criteria for the Design Cells; Diagram of Spaghetti; Preparation of the Work Force; Elements of Production Control; Kanban; Heijunka Box; Kanbans sizing; Aspects of Deployment Management; Visual management. The Six Sigma Strategy and DMAIC method.'; C = 'The Lean Production Approach; History; Lean Thonking; Waste Categories; DMAIC methodology applied to the Lean Six Sigma; Organizational culture for the Lean; Understanding the Value Stream; Value Stream Mapping; Flow of Material and Information; Characteristics of the Lean Value Stream; Concept of Takt Time; Capability Analysis; Production System Pulled; Creation of continuous flow of production; Manufacturing Cells; Criteria for the Design Cells; Diagram of Spaghetti; Preparation of the Work Force; Elements of Production Control; Kanban; Heijunka Box; Kanbans sizing; Aspects of Deployment Management; Visual management. The Six Sigma Strategy and DMAIC method.' },
    @{ Row = 17; A = 'Avaliação:'; B = $null; C = $null },
    @{ Row = 18; A = 'Método:'; B = '5840535 - Messias Borges Silva'; C = '5840535 - Messias Borges Silva' },
    @{ Row = 19; A = 'Critério:'; B = 'Aulas Expositivas; trabalhos e seminários.'; C = 'Aulas Expositivas; trabalhos e seminários.' },
    @{ Row = 20; A = 'Norma de recuperação:'; B = 'MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.'; C = 'MF = (0,30*P1 + 0,30*P2 + 0,40*TRAB), onde P1 e P2 são provas e TRAB é a nota média de trabalhos e seminários.' },
    @{ Row = 21; A = 'Bibliografia:'; B = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação'; C = 'NF = (MF + PR)/2, onde PR é uma prova de recuperação' },
    @{ Row = 22; A = 'Requisitos:'; B = $null; C = $null },
    @{ Row = 23; A = $null; B = "LOQ4260 -  Controle Estatístico da Qualidade  (Requisito fraco)`n"; C = "LOQ4260 -  Controle Estatístico da Qualidade  (Requisito fraco)`n" }
)

foreach ($entry in $rowData) {
    $r = $entry.Row
    if ($null -ne $entry.A) { $ws.Cells.Item($r, 1).Value = $entry.A }
    if ($null -ne $entry.B) { $ws.Cells.Item($r, 2).Value = $entry.B }
    if ($null -ne $entry.C) { $ws.Cells.Item($r, 3).Value = $entry.C }
}

# A few of the rows above used to be completely empty in columns B/C, so the
# new cells inherited the wrong column style from COM's "nearest column"
# default. Re-stamp their number/format (but not the text we just set) by
# copying formats from a same-column neighbour that already has the correct
# style, then restore the value.
function Copy-CellFormat($fromRow, $toRow, $col) {
    $ws.Cells.Item($fromRow, $col).Copy()
    $ws.Cells.Item($toRow, $col).PasteSpecial(-4122)  # xlPasteFormats
}

Copy-CellFormat -fromRow 19 -toRow 18 -col 2
$ws.Cells.Item(18, 2).Value = '5840535 - Messias Borges Silva'
Copy-CellFormat -fromRow 19 -toRow 18 -col 3
$ws.Cells.Item(18, 3).Value = '5840535 - Messias Borges Silva'

Copy-CellFormat -fromRow 20 -toRow 23 -col 2
$ws.Cells.Item(23, 2).Value = "LOQ4260 -  Controle Estatístico da Qualidade  (Requisito fraco)`n"
Copy-CellFormat -fromRow 20 -toRow 23 -col 3
$ws.Cells.Item(23, 3).Value = "LOQ4260 -  Controle Estatístico da Qualidade  (Requisito fraco)`n"

# Row-height changes between the old and new layout.
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30

# The sheet now ends at row 23 — drop the old trailing row 24.
$ws.Rows.Item(24).Delete()
